$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing B94:B102 values (revised UCSV estimates) ---
$updatedValues = @{
    94  = 2.122782707214355
    95  = 2.245511293411255
    96  = 2.197686433792114
    97  = 2.24852728843689
    98  = 2.054852485656738
    99  = 2.121228933334351
    100 = 2.344321250915527
    101 = 2.253693580627441
    102 = 2.123221397399902
}

foreach ($row in $updatedValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $updatedValues[$row]
}

# --- Append new rows 103:114 with new dates + values ---
$newRows = @(
    @{ Row = 103; Year = 2020; Month = 4;  Day = 1; Value = 1.176671504974365 },
    @{ Row = 104; Year = 2020; Month = 7;  Day = 1; Value = 1.706722378730774 },
    @{ Row = 105; Year = 2020; Month = 10; Day = 1; Value = 1.633168458938599 },
    @{ Row = 106; Year = 2021; Month = 1;  Day = 1; Value = 1.660865187644958 },
    @{ Row = 107; Year = 2021; Month = 4;  Day = 1; Value = 4.408313751220703 },
    @{ Row = 108; Year = 2021; Month = 7;  Day = 1; Value = 4.022122859954834 },
    @{ Row = 109; Year = 2021; Month = 10; Day = 1; Value = 5.522685050964355 },
    @{ Row = 110; Year = 2022; Month = 1;  Day = 1; Value = 6.452416896820068 },
    @{ Row = 111; Year = 2022; Month = 4;  Day = 1; Value = 5.884917736053467 },
    @{ Row = 112; Year = 2022; Month = 7;  Day = 1; Value = 6.64294958114624 },
    @{ Row = 113; Year = 2022; Month = 10; Day = 1; Value = 5.703855991363525 },
    @{ Row = 114; Year = 2023; Month = 1;  Day = 1; Value = 5.602568626403809 }
)

# Reference style cell for the date column (row 102, column A) that already
# carries the desired number format / font / border.
$styleSource = $ws.Cells.Item(102, 1)

foreach ($entry in $newRows) {
    $row = $entry.Row
    $dateCell = $ws.Cells.Item($row, 1)
    $valueCell = $ws.Cells.Item($row, 2)

    $styleSource.Copy()
    $dateCell.PasteSpecial(-4122)

    $d = Get-Date -Year $entry.Year -Month $entry.Month -Day $entry.Day
    $dateCell.Value = $d.Date
    $valueCell.Value = $entry.Value
}
